$wb = $excel.ActiveWorkbook

# Rename the translation attributes sheet: "TRN-1079-0833-9890 (FA)" -> "FA (TRN-1079-0833-9890)"
$ws = $wb.Worksheets.Item("TRN-1079-0833-9890 (FA)")
$ws.Name = "FA (TRN-1079-0833-9890)"

# Fix the header row: the old export swapped the "Key"/"Original Value"/"Action"/"Value"/"Comment"
# headers; the corrected export uses lowercase column keys plus a new "editor" column, and the
# "Original Value" header becomes the locale-specific label.
$ws.Range("A1").Value2 = "key"
$ws.Range("B1").Value2 = "Persian (TRN-1079-0833-9890)"
$ws.Range("C1").Value2 = "action"
$ws.Range("D1").Value2 = "value"
$ws.Range("E1").Value2 = "comment"

# Add the new "editor" column header, matching the formatting of the other plain headers.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value2 = "editor"

# Populate the new "editor" column for every data row with the same placeholder used by
# the "action" column.
for ($r = 2; $r -le 31; $r++) {
    $cell = "F" + $r
    $ws.Range($cell).Value2 = "-"
}

# Restore the active selection on the sheet.
$ws.Activate()
$ws.Range("E18").Select()
